# refactor: rewrite encoder (#5)
#
# The "survey" sheet's header cell C1 held the label "label:English";
# it is corrected to the XLSForm-style "label::English (en)" (matching
# the same text already used as the header of the "choices" sheet).
# The "survey" tab becomes the active/selected sheet (it previously was
# "settings"), with the cursor left on C2 just below the edited header.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")

# Fix the mislabeled survey header (matches the text already used on the
# "choices" sheet's own label column header).
$survey.Range("C1").Value = "label::English (en)"

# Make "survey" the active sheet/tab (it was "settings") and leave the
# selection on C2, mirroring the cursor move after editing the header cell.
$survey.Activate()
$survey.Range("C2").Select()
